$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 14.333333
$ws.Range("I8").Value = 14.333333
$ws.Range("K8").Value = 42.999999
$ws.Range("M8").Value = 96.000001

$ws.Range("H33").Value = 257.91666
$ws.Range("I33").Value = 229.5
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 229.5
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = -0.5
$ws.Range("N33").Value = -858

$ws.Range("H51").Value = 53870.477
$ws.Range("I51").Value = 9170.571
$ws.Range("J51").Value = 76220.42999999999
$ws.Range("K51").Value = 9170.571
$ws.Range("L51").Value = 76220.42999999999
$ws.Range("M51").Value = -8686.571
$ws.Range("N51").Value = -77188.42999999999

$ws.Range("H74").Value = 3998.5
$ws.Range("I74").Value = 3998.5
$ws.Range("K74").Value = 3998.5
$ws.Range("M74").Value = -3062.5

$ws.Range("H76").Value = 4000
$ws.Range("J76").Value = 4000
$ws.Range("L76").Value = 4000
$ws.Range("N76").Value = -4630

$ws.Range("H77").Value = 3998.5
$ws.Range("I77").Value = 3998.5
$ws.Range("K77").Value = 19992.5
$ws.Range("M77").Value = -15312.5

$ws.Range("H79").Value = 4000
$ws.Range("J79").Value = 4000
$ws.Range("L79").Value = 4000
$ws.Range("N79").Value = -6184

$ws.Range("H99").Value = 590
$ws.Range("J99").Value = 590
$ws.Range("L99").Value = 1770
$ws.Range("N99").Value = -4766

$ws.Range("H121").Value = 5776.6
$ws.Range("J121").Value = 5776.6
$ws.Range("L121").Value = 17329.8
$ws.Range("N121").Value = -20823.8

$ws.Range("H129").Value = 2163.889
$ws.Range("I129").Value = 1338.2
$ws.Range("J129").Value = 3196
$ws.Range("K129").Value = 4014.6
$ws.Range("L129").Value = 9588
$ws.Range("M129").Value = 985.3999999999996
$ws.Range("N129").Value = -19588

$ws.Range("H137").Value = 931.93335
$ws.Range("I137").Value = 931.93335
$ws.Range("K137").Value = 2795.80005
$ws.Range("M137").Value = -245.8000499999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1698.75
$ws.Range("I5").Value = 1598.3334
$ws.Range("K5").Value = 1598.3334
$ws.Range("M5").Value = -1486.3334

$ws.Range("H11").Value = 17500638
$ws.Range("I11").Value = 17500638
$ws.Range("K11").Value = 17500638
$ws.Range("M11").Value = -17500494

$ws.Range("H16").Value = 5481.2
$ws.Range("I16").Value = 6101.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 6101.5
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -5814.5
$ws.Range("N16").Value = -3574

$ws.Range("H17").Value = 8400
$ws.Range("J17").Value = 8400
$ws.Range("L17").Value = 8400
$ws.Range("N17").Value = -8746

$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 2000
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -3372

$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 10000
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -16864

$ws.Range("H102").Value = 74074744
$ws.Range("I102").Value = 74074744
$ws.Range("K102").Value = 74074744
$ws.Range("M102").Value = -74073122

$ws.Range("H110").Value = 3833701.8
$ws.Range("I110").Value = 6174330.5
$ws.Range("K110").Value = 6174330.5
$ws.Range("M110").Value = -6172285.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1698.75
$ws.Range("I4").Value = 1598.3334
$ws.Range("K4").Value = 1598.3334
$ws.Range("M4").Value = -1483.3334

$ws.Range("H64").Value = 1332.8334
$ws.Range("I64").Value = 1166.6666
$ws.Range("J64").Value = 1499
$ws.Range("K64").Value = 1166.6666
$ws.Range("L64").Value = 1499
$ws.Range("M64").Value = -941.6666
$ws.Range("N64").Value = -1949

$ws.Range("H67").Value = 1332.8334
$ws.Range("I67").Value = 1166.6666
$ws.Range("J67").Value = 1499
$ws.Range("K67").Value = 1166.6666
$ws.Range("L67").Value = 1499
$ws.Range("M67").Value = -386.6666
$ws.Range("N67").Value = -3059

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 17000
$ws.Range("I26").Value = 14500
$ws.Range("J26").Value = 22000
$ws.Range("K26").Value = 14500
$ws.Range("L26").Value = 22000
$ws.Range("M26").Value = -14213
$ws.Range("N26").Value = -22574

$ws.Range("H28").Value = 15799.6
$ws.Range("J28").Value = 17249.5
$ws.Range("L28").Value = 17249.5
$ws.Range("N28").Value = -17739.5

$ws.Range("H56").Value = 5062
$ws.Range("I56").Value = 5093
$ws.Range("K56").Value = 5093
$ws.Range("M56").Value = -4248

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

$ws.Range("H68").Value = 494.7143
$ws.Range("I68").Value = 391
$ws.Range("J68").Value = 633
$ws.Range("K68").Value = 1173
$ws.Range("L68").Value = 1899
$ws.Range("M68").Value = -362
$ws.Range("N68").Value = -3521

$ws.Range("H71").Value = 494.7143
$ws.Range("I71").Value = 391
$ws.Range("J71").Value = 633
$ws.Range("K71").Value = 3519
$ws.Range("L71").Value = 5697
$ws.Range("M71").Value = 537
$ws.Range("N71").Value = -13809

$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 12000
$ws.Range("N80").Value = -13872

$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 36000
$ws.Range("N83").Value = -45360

$ws.Range("H109").Value = 1537.2858
$ws.Range("J109").Value = 2996.5
$ws.Range("L109").Value = 8989.5
$ws.Range("N109").Value = -11069.5

$ws.Range("H130").Value = 3468.75
$ws.Range("I130").Value = 2315
$ws.Range("J130").Value = 4622.5
$ws.Range("K130").Value = 6945
$ws.Range("L130").Value = 13867.5
$ws.Range("M130").Value = -1925
$ws.Range("N130").Value = -23907.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2916.6667
$ws.Range("I31").Value = 2916.6667
$ws.Range("K31").Value = 2916.6667
$ws.Range("M31").Value = -2624.6667

$ws.Range("H37").Value = 2916.6667
$ws.Range("I37").Value = 2916.6667
$ws.Range("K37").Value = 2916.6667
$ws.Range("M37").Value = -2639.6667

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H70").Value = 14758.2
$ws.Range("I70").Value = 7800
$ws.Range("J70").Value = 16497.75
$ws.Range("K70").Value = 7800
$ws.Range("L70").Value = 16497.75
$ws.Range("M70").Value = -7530
$ws.Range("N70").Value = -17037.75

$ws.Range("H73").Value = 14758.2
$ws.Range("I73").Value = 7800
$ws.Range("J73").Value = 16497.75
$ws.Range("K73").Value = 7800
$ws.Range("L73").Value = 16497.75
$ws.Range("M73").Value = -6864
$ws.Range("N73").Value = -18369.75

$ws.Range("H97").Value = 1468.0416
$ws.Range("I97").Value = 1109.6364
$ws.Range("J97").Value = 1771.3077
$ws.Range("K97").Value = 1109.6364
$ws.Range("L97").Value = 1771.3077
$ws.Range("M97").Value = -613.6364000000001
$ws.Range("N97").Value = -2763.3077

$ws.Range("H126").Value = 6734.5
$ws.Range("I126").Value = 6181.4
$ws.Range("K126").Value = 18544.2
$ws.Range("M126").Value = -16074.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4163.2856
$ws.Range("I7").Value = 3329.6924
$ws.Range("K7").Value = 3329.6924
$ws.Range("M7").Value = -3217.6924

$ws.Range("H22").Value = 17909
$ws.Range("I22").Value = 19000
$ws.Range("J22").Value = 17499.875
$ws.Range("K22").Value = 19000
$ws.Range("L22").Value = 17499.875
$ws.Range("M22").Value = -18705
$ws.Range("N22").Value = -18089.875

$ws.Range("H24").Value = 40007
$ws.Range("J24").Value = 40007
$ws.Range("L24").Value = 40007
$ws.Range("N24").Value = -40693

$ws.Range("H27").Value = 17909
$ws.Range("I27").Value = 19000
$ws.Range("J27").Value = 17499.875
$ws.Range("K27").Value = 19000
$ws.Range("L27").Value = 17499.875
$ws.Range("M27").Value = -18893
$ws.Range("N27").Value = -17713.875

$ws.Range("H68").Value = 5500
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 5500
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 5500
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -6998

$ws.Range("H71").Value = 5500
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 5500
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 27500
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -34988

$ws.Range("H82").Value = 78389
$ws.Range("I82").Value = 2666
$ws.Range("J82").Value = 101105.9
$ws.Range("K82").Value = 2666
$ws.Range("L82").Value = 101105.9
$ws.Range("M82").Value = -2305
$ws.Range("N82").Value = -101827.9

$ws.Range("H85").Value = 78389
$ws.Range("I85").Value = 2666
$ws.Range("J85").Value = 101105.9
$ws.Range("K85").Value = 2666
$ws.Range("L85").Value = 101105.9
$ws.Range("M85").Value = -1418
$ws.Range("N85").Value = -103601.9

$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()

$ws.Range("H126").Value = 4163.2856
$ws.Range("I126").Value = 3329.6924
$ws.Range("K126").Value = 9989.0772
$ws.Range("M126").Value = -7519.0772

$ws.Range("H132").Value = 145062.28
$ws.Range("I132").Value = 168869
$ws.Range("K132").Value = 506607
$ws.Range("M132").Value = -504077

$ws.Range("H136").Value = 5761.125
$ws.Range("I136").Value = 4500.7
$ws.Range("J136").Value = 7861.8335
$ws.Range("K136").Value = 13502.1
$ws.Range("L136").Value = 23585.5005
$ws.Range("M136").Value = -10952.1
$ws.Range("N136").Value = -28685.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 417400
$ws.Range("J26").Value = 417400
$ws.Range("L26").Value = 417400
$ws.Range("N26").Value = -417986
